{"js": "// The report body was substantially rewritten: the narrative paragraphs were\n// reworded, the flat \"Pending Updates\" list was replaced with a package/arch\n// list plus a review note, the RMF-compliance / next-steps / risk-assessment\n// sections were rewritten from numbered lists into prose + \"- \" bullet lists,\n// and the closing paragraphs were replaced. The paragraph count does not\n// change (34 before and after), so every paragraph can be addressed by its\n// (stable) index and have its text replaced in place.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// index -> new paragraph text (\\u000b == the manual line break rendered as <w:br/>)\nconst newParagraphText = {\n  2: \"The network is composed of several systems with different update structures. The code/stable and ure/stable-security repositories contain patches for the most recent versions, while the git-man/stable-security and git/stable-security repositories have more stable versions.\\u000b\",\n  5: \"There are pending updates available for the following systems: \\u000b\",\n  6: \"\\u000b\",\n  7: \"- code/stable 1.99.0-1743632463 amd64\\u000b\",\n  8: \"- ure/stable-security 4:7.4.7-1+deb12u6 amd64 \\u000b\",\n  9: \"- git-man/stable-security 1:2.39.5-0+deb12u2 all \\u000b\",\n  10: \"- git/stable-security 1:2.39.5-0+deb12u2 amd64 \\u000b\",\n  11: \"\\u000b\",\n  12: \"These updates are currently in the process of being reviewed and verified.\\u000b\",\n  13: \"\\u000b\",\n  14: \"*** Compliance with RMF Controls ***\\u000b\",\n  15: \"To ensure compliance with the Risk Management Framework (RMF), it is essential to identify, report, and take corrective action on any vulnerabilities found. The following steps should be taken:\\u000b\",\n  16: \"\\u000b\",\n  17: \"- Identification: Identify the affected systems and their corresponding patch versions.\\u000b\",\n  18: \"- Reporting: Document the identified vulnerabilities, including the potential impact level and mitigation plan.\\u000b\",\n  19: \"- Corrective Action: Apply the necessary patches to remediate the vulnerabilities.\\u000b\",\n  21: \"In this scenario, the identified vulnerabilities include Command Injection in certain Git repositories and a lack of certificate validation in CODESYS Git. Proper configuration management and vulnerability checks should be implemented to prevent similar incidents in the future.\\u000b\",\n  22: \"\\u000b\",\n  23: \"*** Recommended next steps ***\\u000b\",\n  24: \"The following actions are recommended to ensure the successful implementation of the patch:\\u000b\",\n  26: \"- Review and assess the updates available for all systems.\\u000b\",\n  27: \"- Schedule the deployment of patches for each system, ensuring that no system is left without an update.\\u000b\",\n  28: \"- Update documentation to reflect any changes or modifications made during the patching process.\\u000b\",\n  29: \"\\u000b\",\n  30: \"*** Risk Assessment ***\\u000b\",\n  31: \"The network is at risk due to several identified vulnerabilities. The potential impact level is moderate to high, depending on the severity and extent of the vulnerability. The mitigation plan includes applying necessary patches, configuring systems for proper security, and implementing vulnerability checks. Regular monitoring and review are essential to ensure the effectiveness of these measures.\\u000b\",\n  33: \"In this scenario, the potential risk is due to the lack of certificate validation in CODESYS Git and other vulnerabilities found in various repositories. The impact level is moderate to high, as unauthorized access or manipulation of sensitive data could occur if not addressed promptly.\",\n};\n\nfor (const [indexStr, text] of Object.entries(newParagraphText)) {\n  const index = Number(indexStr);\n  const paragraph = paragraphs.items[index];\n  if (paragraph.text !== text) {\n    // \"Replace\" rewrites this paragraph's run text (and its trailing <w:br/>,\n    // driven by the \\u000b) without touching neighboring paragraphs.\n    paragraph.insertText(text, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "# The report body was substantially rewritten: the narrative paragraphs were\n# reworded, the flat \"Pending Updates\" list was replaced with a package/arch\n# list plus a review note, the RMF-compliance / next-steps / risk-assessment\n# sections were rewritten from numbered lists into prose + \"- \" bullet lists,\n# and the closing paragraphs were replaced. The paragraph count does not\n# change (34 before and after), so every paragraph can be addressed by its\n# (stable, 1-based) Paragraphs() index and have its Range.Text replaced in\n# place. [char]11 is the manual line break (<w:br/>) Word stores as \\v.\n$d = $word.ActiveDocument\n\n$newParagraphText = [ordered]@{\n  3 = 'The network is composed of several systems with different update structures. The code/stable and ure/stable-security repositories contain patches for the most recent versions, while the git-man/stable-security and git/stable-security repositories have more stable versions.' + [char]11\n  6 = 'There are pending updates available for the following systems: ' + [char]11\n  7 = '' + [char]11\n  8 = '- code/stable 1.99.0-1743632463 amd64' + [char]11\n  9 = '- ure/stable-security 4:7.4.7-1+deb12u6 amd64 ' + [char]11\n  10 = '- git-man/stable-security 1:2.39.5-0+deb12u2 all ' + [char]11\n  11 = '- git/stable-security 1:2.39.5-0+deb12u2 amd64 ' + [char]11\n  12 = '' + [char]11\n  13 = 'These updates are currently in the process of being reviewed and verified.' + [char]11\n  14 = '' + [char]11\n  15 = '*** Compliance with RMF Controls ***' + [char]11\n  16 = 'To ensure compliance with the Risk Management Framework (RMF), it is essential to identify, report, and take corrective action on any vulnerabilities found. The following steps should be taken:' + [char]11\n  17 = '' + [char]11\n  18 = '- Identification: Identify the affected systems and their corresponding patch versions.' + [char]11\n  19 = '- Reporting: Document the identified vulnerabilities, including the potential impact level and mitigation plan.' + [char]11\n  20 = '- Corrective Action: Apply the necessary patches to remediate the vulnerabilities.' + [char]11\n  22 = 'In this scenario, the identified vulnerabilities include Command Injection in certain Git repositories and a lack of certificate validation in CODESYS Git. Proper configuration management and vulnerability checks should be implemented to prevent similar incidents in the future.' + [char]11\n  23 = '' + [char]11\n  24 = '*** Recommended next steps ***' + [char]11\n  25 = 'The following actions are recommended to ensure the successful implementation of the patch:' + [char]11\n  27 = '- Review and assess the updates available for all systems.' + [char]11\n  28 = '- Schedule the deployment of patches for each system, ensuring that no system is left without an update.' + [char]11\n  29 = '- Update documentation to reflect any changes or modifications made during the patching process.' + [char]11\n  30 = '' + [char]11\n  31 = '*** Risk Assessment ***' + [char]11\n  32 = 'The network is at risk due to several identified vulnerabilities. The potential impact level is moderate to high, depending on the severity and extent of the vulnerability. The mitigation plan includes applying necessary patches, configuring systems for proper security, and implementing vulnerability checks. Regular monitoring and review are essential to ensure the effectiveness of these measures.' + [char]11\n  34 = 'In this scenario, the potential risk is due to the lack of certificate validation in CODESYS Git and other vulnerabilities found in various repositories. The impact level is moderate to high, as unauthorized access or manipulation of sensitive data could occur if not addressed promptly.'\n}\n\nforeach ($index in $newParagraphText.Keys) {\n  $d.Paragraphs($index).Range.Text = $newParagraphText[$index]\n}"}
